$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: names; Column B: pixel values
$names = @(
  "name",
  "Third_eyelid_protrude",
  "lens_escape",
  "epiphora_brown",
  "blepharitis_inflammation",
  "blepharitis_inner_inflammation",
  "nuclear_sclerosis_gray",
  "glaucoma_macrophthalmia",
  "glaucoma_serosity",
  "glaucoma_blue_white",
  "glaucoma_flare",
  "glaucoma_anisocoria",
  "glaucoma_lens_escape",
  "keratitis_Lipid_deposit",
  "keratitis_pigmentation",
  "keratitis_flare",
  "keratoconjunctivitis_sicca_mucus",
  "keratoconjunctivitis_sicca_pigmentation",
  "keratoconjunctivitis_sicca_flare",
  "keratoconjunctivitis_sicca_corneal_opacity",
  "uveitis_discoloration",
  "uveitis_miosis",
  "uveitis_flare",
  "uveitis_inflammation",
  "uveitis_aqueous_flare",
  "trichiasis_hair",
  "corneal_pus",
  "corneal_scratch",
  "corneal",
  "conjunctivitis_flare",
  "conjunctivitis_swll",
  "conjunctivitis_white_inflammation",
  "gataract",
  "gataract_initial",
  "ectropion_droop",
  "ectropion_roll_in",
  "ectropion_damage",
  "ectropion_corneal_damage",
  "corneal_edema_moon_halo",
  "corneal_edema_irradiation",
  "corneal_degeneration_Glass_fiber",
  "corneal_degeneration_opacity"
)

for ($i = 0; $i -lt $names.Count; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $names[$i]
  if ($row -eq 1) {
    $ws.Cells.Item($row, 2).Value = "pixelValue"
  } else {
    $ws.Cells.Item($row, 2).Value = $row - 1
  }
}

# New rows 41-42 (corneal_degeneration_*) are vertically centered
$ws.Range("A41:B42").VerticalAlignment = -4108

# Column A width
$ws.Columns("A").ColumnWidth = 35.25

# View: zoomed out to 70%, selection moved to A20
$excel.ActiveWindow.Zoom = 70
[void]$ws.Range("A20").Select()
